$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Step 1: create the 10 new shared strings (indices 34-43) in the exact order
# required, by writing them into their destination cells in index order --
# Excel appends newly-seen strings to the shared-string table in first-use order.
$ws.Range("C10").Value2 = ' I\''ve got a bad feeling about the\n[CS:N]Sneasel[CR] at the other table for some reason.'
$ws.Range("C11").Value2 = ' As soon as he came into the\nstore, he yelled at me, \"Don\''t just loiter!\"'
$ws.Range("C12").Value2 = ' It\''s OK to take your time and\nbrowse a bit, though.[K] Isn\''t it?'
$ws.Range("A10").Value2 = 'SCRIPT/P01P04A/um1103.ssb'
$ws.Range("D10").Value2 = ' Я не знаю почему, но мне не\nнравится [CS:N]Снизел[CR] за соседним столиком.'
$ws.Range("D12").Value2 = ' Но в том, чтобы побродить\nтуда-сюда нет ничего такого.[K] Так ведь?'
$ws.Range("E10").Value2 = ' Ÿ îå èîàý ðïœåíô, îï íîå îå\nîñàâéóòÿ [CS:N]Òîéèåì[CR] èà òïòåäîéí òóïìéëïí.'
$ws.Range("E11").Value2 = ' Ïî èàšæì òýäà é óôó çå îàœàì\nîà íåîÿ ëñéœàóû: \"Öâàóéó òìïîÿóûòÿ áåè\näåìà!\"'
$ws.Range("E12").Value2 = ' Îï â óïí, œóïáú ðïáñïäéóû\nóôäà-òýäà îåó îéœåãï óàëïãï.[K] Óàë âåäû?'
$ws.Range("D11").Value2 = ' Он зашёл сюда и тут же начал\nна меня кричать: \"Хватит слоняться без\nдела!\"'

# Step 2: numeric "line number" column values for the new rows
$ws.Range("B10").Value2 = 258
$ws.Range("B11").Value2 = 267
$ws.Range("B12").Value2 = 270

# Step 3: row 12 uses the "last row" border style (like existing row 9) --
# copy formatting only from row 9 onto row 12 (also sets A12 style with no value).
$ws.Range("A9:E9").Copy()
$ws.Range("A12:E12").PasteSpecial(-4122)
[void]$excel.CutCopyMode

# Step 4: row heights matching the wrapped content (same heights as used
# elsewhere in the sheet for similarly-sized text).
$ws.Rows.Item(10).RowHeight = 43.2
$ws.Rows.Item(11).RowHeight = 31.8
$ws.Rows.Item(12).RowHeight = 21.6

# Step 5: update the view selection to match the author's final cursor position.
[void]$ws.Range("B7").Select()
